# Actualización automática hashcode sáb oct  5 01:55:50 CEST 2019
# Updates the MD5 hashcode values (column B) for a set of rows identified
# by their key in column A, matching the upstream CSV->xlsx re-export diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B9";   Old = "b38f934c02d047a2ada11101a82c1f39"; New = "07bacf8e347e2ede09f2dd7b7ce23ccd" },
    @{ Cell = "B17";  Old = "07256692167359f375548b4159378639"; New = "439a4ecb56e117e8046f17ba1bf3a6f7" },
    @{ Cell = "B34";  Old = "82a122538dd440102d3a80b6a21db178"; New = "00198639910968560f11c8d8bade01df" },
    @{ Cell = "B94";  Old = "44213aeeab26b84a909d27da8747f1dd"; New = "3976bbb9f1b4a382bc87fb541bc59088" },
    @{ Cell = "B95";  Old = "62fb3a25e5eb73fa548e78df049eeae4"; New = "5c73882c3c53b385a9b6cb3418168164" },
    @{ Cell = "B98";  Old = "7e28e709da59e3fc566edfc13a487028"; New = "dfb77a4ff63d5cca57d0b52f6e8ac4ad" },
    @{ Cell = "B109"; Old = "4eadddab98df18409f53e51a7d916afb"; New = "be3773578f4667dab4f496d33a85eeab" },
    @{ Cell = "B115"; Old = "78fb34603fc974bb8815be6ff28d67f3"; New = "662197525b2acd21c6124243032fe1bd" },
    @{ Cell = "B162"; Old = "496da3c040126f1aa643fcc0bd0ac7b4"; New = "f51c5801efcb5a7779994e0ab01ab31c" },
    @{ Cell = "B175"; Old = "a0415eaf12e05dffcf13385a258b6d22"; New = "341d04dc8624b06de520ab052b246705" },
    @{ Cell = "B180"; Old = "3628b7505f9fe43df36ba6974d4ef11f"; New = "ea64f177d870ab1b4d17195e11563195" },
    @{ Cell = "B183"; Old = "b3a5b41de62bc70708855999dc05272a"; New = "5604193cd3e0c50da54a7cc418e8c0d1" },
    @{ Cell = "B200"; Old = "167b2fa8a52251f81750b9c2cb5d4eea"; New = "54fde3c652374a90a060cb0833078d4b" },
    @{ Cell = "B213"; Old = "6d257cf1531177a2c618d10a50546c6c"; New = "bc845d51fe4526a35fa2a6a0a2d9f2ff" },
    @{ Cell = "B227"; Old = "2d01a5278488f10b9f5dd5e43c9859b6"; New = "0bed25d524905a11baaf024e5fd8abc9" },
    @{ Cell = "B228"; Old = "64b0b49079d4fafbf463562b0ce5c243"; New = "e78f6065457a7f20686dac6a2bde44af" },
    @{ Cell = "B232"; Old = "c7017acfe56676dd01830aabf3c16619"; New = "72e804d3ceaaf08953cc162b25b3431f" },
    @{ Cell = "B302"; Old = "128c4596fca9a98de68b10dcf6d5b902"; New = "eb82f80ba4867ba6d0317941901c73a3" },
    @{ Cell = "B420"; Old = "bf3569543f5afe0bd329968445d710df"; New = "0841f66eec1f7caf51680bed6f5054c6" },
    @{ Cell = "B465"; Old = "89c67370eabfd551687d12306ce287f7"; New = "0ba828f890635471b5f48b686a7528ff" },
    @{ Cell = "B483"; Old = "7db025c699f5ae5fc290487270fbbc2d"; New = "0a8277e209a3872254017c3c03014112" },
    @{ Cell = "B513"; Old = "7ae2c5bb5dacbf5ba8bf260171240429"; New = "41eff1cb20baf8ee2e96f81bad43e7b4" },
    @{ Cell = "B558"; Old = "48430e2174399aad2d97e1908c082c03"; New = "4d73a6bc8ae7518487d96671ff690408" },
    @{ Cell = "B580"; Old = "fa0233183a94dd823d1a0c00a9af25d2"; New = "e3c0ec123128990d05c2873928241288" },
    @{ Cell = "B600"; Old = "98a7a4c7e45a4c7f13b04e8c8f695464"; New = "eef16b95de2fdd043b7a987a50adf02f" },
    @{ Cell = "B626"; Old = "cdeec3a4e361cc7e3e633c7a2be1280d"; New = "124054d4a6a4cbe2c5a28c761a12800c" },
    @{ Cell = "B673"; Old = "101774f99322477ca9526553e92d1dd0"; New = "b8c3365c89986858fe3f978fef075e67" },
    @{ Cell = "B723"; Old = "356ca7a6a0143f6e4c614d0549b08df8"; New = "f42b4c242b8194372d36b6b3b83920e7" },
    @{ Cell = "B737"; Old = "8ab5bc0100be605a7e27d1c9c2d71284"; New = "6110daed30891fc642307af5566e649a" },
    @{ Cell = "B862"; Old = "cabac408ee7be64c2ee1efcd01eb2d8a"; New = "e32da4b4d818328c59a076bd5bd27191" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $current = $cell.Value2
    if ($current -ne $u.Old) {
        Write-Output "Warning: $($u.Cell) expected '$($u.Old)' but found '$current'"
    }
    $cell.Value = $u.New
}
